$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F10").Value = "WIP"
$ws.Range("G10").Value = "1. Basics-VII.py (Base and original assignment)"
